# windows iRave problem and excell update
# Re-creates the data-entry updates made to labs/lab13/excel.xlsx:
#  - corrects Tarefa 2 / Tarefa 3 "SEQ" scores for participant 1 (row 8)
#  - fixes a typo in participant 6's occupation (row 13)
#  - fills in the previously-empty rows for participants 7-9 (rows 14-16),
#    including a new "age bracket" column (U)
#  - updates the current selection/active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8: Tarefa 2 / Tarefa 3 SEQ values corrected ---
$ws.Range("F8").Value = 4
$ws.Range("J8").Value = 5
$ws.Range("N8").Value = 6

# --- Row 14 (participant 7) ---
$ws.Range("D14").Value = 28
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 3
$ws.Range("H14").Value = 47
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 4
$ws.Range("L14").Value = 60
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 5
$ws.Range("P14").Value = "63 anos"
$ws.Range("Q14").Value = "M"
$ws.Range("R14").Value = "Professor"
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 4
$ws.Range("U14").Value = "0-10"

# --- Row 15 (participant 8) ---
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 22
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 6
$ws.Range("L15").Value = 45
$ws.Range("M15").Value = 2
$ws.Range("N15").Value = 4
$ws.Range("P15").Value = "16 anos"
$ws.Range("Q15").Value = "F"
$ws.Range("R15").Value = "estudante"
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 6
$ws.Range("U15").Value = "0-10"

# --- Row 16 (participant 9) ---
$ws.Range("D16").Value = 17
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = 59
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("L16").Value = 43
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 5
$ws.Range("P16").Value = "14 anos"
$ws.Range("Q16").Value = "F"
$ws.Range("R16").Value = "estudante"
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = 8
$ws.Range("U16").Value = "0-10"

# --- Row 13: fix occupation typo ("Cabelareira" -> "cabeleireiro") ---
$ws.Range("R13").Value = "cabeleireiro"

# --- Update the active selection to match the saved view ---
$ws.Activate()
$ws.Range("N8").Select()
